$d = $word.ActiveDocument

# The first paragraph currently reads:
#   "Screenshot of MongoDB Compass:  " + "mars_app"
# split across two runs with spell-check proofErr markers around "mars_app".
# Round-trip that paragraph's Range through WordOpenXML/InsertXML so Word
# collapses it into a single clean run with no leftover proofErr markers.
$p1 = $d.Paragraphs.Item(1)
$cleanXml = $p1.Range.WordOpenXML
[void]$p1.Range.InsertXML($cleanXml)

# Now insert a brand new paragraph above it containing the new heading text.
$p1 = $d.Paragraphs.Item(1)
$p1.Range.InsertParagraphBefore()
$titlePara = $d.Paragraphs.Item(1)
$titlePara.Range.Text = "Web Scraping Screenshots"
